$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fully clear A1 (value + formatting) since header moved from A:E to B:J and A1 is now blank
$ws.Range("A1").Clear()

# New header row (B1:J1)
$ws.Range("B1").Value = "CAGE#"
$ws.Range("C1").Value = "Gene"
$ws.Range("D1").Value = "Guide"
$ws.Range("E1").Value = "Replicate"
$ws.Range("F1").Value = "init_oof"
$ws.Range("G1").Value = "final_oof"
$ws.Range("H1").Value = "fitness_score"
$ws.Range("I1").Value = "avg_fit_score"
$ws.Range("J1").Value = "stdev"

# New row-index column A, values 1-4
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4

# Row 2 (CAGE3267 / hFGFR3 / g46 / rep1)
$ws.Range("B2").Value = "CAGE3267"
$ws.Range("C2").Value = "hFGFR3"
$ws.Range("D2").Value = "g46"
$ws.Range("E2").Value = "CAGE3267_RT112-g46_rep1"
$ws.Range("F2").Value = 59.2
$ws.Range("G2").Value = 58
$ws.Range("H2").Value = 0.98
$ws.Range("I2").Value = 0.865
$ws.Range("J2").Value = 0.1626345596729059

# Row 3 (CAGE3267 / hFGFR3 / g46 / rep2)
$ws.Range("B3").Value = "CAGE3267"
$ws.Range("C3").Value = "hFGFR3"
$ws.Range("D3").Value = "g46"
$ws.Range("E3").Value = "CAGE3267_RT112-g46_rep2"
$ws.Range("F3").Value = 70
$ws.Range("G3").Value = 52.5
$ws.Range("H3").Value = 0.75
$ws.Range("I3").Value = 0.865
$ws.Range("J3").Value = 0.1626345596729059

# Row 4 (CAGE3269 / hFGFR3 / g9 / rep1)
$ws.Range("B4").Value = "CAGE3269"
$ws.Range("C4").Value = "hFGFR3"
$ws.Range("D4").Value = "g9"
$ws.Range("E4").Value = "CAGE3269_RT112-g9_rep1"
$ws.Range("F4").Value = 4.3
$ws.Range("G4").Value = 3
$ws.Range("H4").Value = 0.7
$ws.Range("I4").Value = 0.74
$ws.Range("J4").Value = 0.05656854249492385

# Row 5 (CAGE3269 / hFGFR3 / g9 / rep2)
$ws.Range("B5").Value = "CAGE3269"
$ws.Range("C5").Value = "hFGFR3"
$ws.Range("D5").Value = "g9"
$ws.Range("E5").Value = "CAGE3269_RT112-g9_rep2"
$ws.Range("F5").Value = 5.4
$ws.Range("G5").Value = 4.2
$ws.Range("H5").Value = 0.78
$ws.Range("I5").Value = 0.74
$ws.Range("J5").Value = 0.05656854249492385

# Apply the header/bold/border style (same style already used by B1:E1) to the
# newly added header cells F1:J1 and to the new index column A2:A5
$ws.Range("B1").Copy()
$ws.Range("F1:J1").PasteSpecial(-4122)
$ws.Range("A2:A5").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
